$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab a template cell that already carries the "quote prefix" cell style
# (s="1" -> quotePrefix="1" in cellXfs) so we can stamp that same format
# onto the new rows that need it, without disturbing their numeric value.
$styleSource = $ws.Range("A9")
$styleSource.Copy()

# New data for rows 7-22 (A: id_subsector, B: id_unit_user_type)
$newValues = @(
    @(31, 6, $true),
    @(32, 6, $true),
    @(33, 6, $true),
    @(34, 6, $true),
    @(35, 6, $true),
    @(36, 6, $true),
    @(37, 6, $true),
    @(38, 6, $true),
    @(39, 6, $true),
    @(310, 6, $true),
    @(311, 6, $true),
    @(312, 6, $true),
    @(313, 6, $false),
    @(314, 6, $false),
    @(315, 6, $false),
    @(316, 6, $false)
)

$r = 7
foreach ($pair in $newValues) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $pair[0]
    if ($pair[2]) {
        $cellA.PasteSpecial(-4122)
    }
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

$excel.CutCopyMode = $false

# Delete the two now-unused rows at the bottom (23 and 24)
$ws.Range("A23:B24").EntireRow.Delete()

# Update selection to match recorded author action
$ws.Range("B13").Select()
